$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the player stats table (rows 2-37); header row 1 is unchanged.
$rows = 36
$cols = 7
$arr = New-Object 'object[,]' $rows,$cols

$arr[0,0] = "Will Simpson"
$arr[0,1] = 1
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 6
$arr[0,5] = 5
$arr[0,6] = 1
$arr[1,0] = "Coby Lovelace"
$arr[1,1] = 2
$arr[1,2] = 2
$arr[1,3] = 0
$arr[1,4] = 5
$arr[1,5] = 3
$arr[1,6] = 2
$arr[2,0] = "Roman Ramirez"
$arr[2,1] = 10
$arr[2,2] = 4
$arr[2,3] = 0
$arr[2,4] = 5
$arr[2,5] = 1
$arr[2,6] = 1
$arr[3,0] = "Jack Massingill"
$arr[3,1] = 2
$arr[3,2] = 3
$arr[3,3] = 0
$arr[3,4] = 5
$arr[3,5] = 5
$arr[3,6] = 0
$arr[4,0] = "Ann Hall"
$arr[4,1] = 0
$arr[4,2] = 2
$arr[4,3] = 0
$arr[4,4] = 4
$arr[4,5] = 3
$arr[4,6] = 1
$arr[5,0] = "Aaron Carter"
$arr[5,1] = 7
$arr[5,2] = 7
$arr[5,3] = 0
$arr[5,4] = 4
$arr[5,5] = 0
$arr[5,6] = 0
$arr[6,0] = "Yvonne Nguyen"
$arr[6,1] = 0
$arr[6,2] = 0
$arr[6,3] = 0
$arr[6,4] = 4
$arr[6,5] = 4
$arr[6,6] = 0
$arr[7,0] = "Gabe Silverstein"
$arr[7,1] = 4
$arr[7,2] = 3
$arr[7,3] = 0
$arr[7,4] = 4
$arr[7,5] = 5
$arr[7,6] = 0
$arr[8,0] = "Jason Jackson"
$arr[8,1] = 5
$arr[8,2] = 4
$arr[8,3] = 1
$arr[8,4] = 3
$arr[8,5] = 3
$arr[8,6] = 1
$arr[9,0] = "Kevin Lee"
$arr[9,1] = 13
$arr[9,2] = 5
$arr[9,3] = 0
$arr[9,4] = 3
$arr[9,5] = 0
$arr[9,6] = 0
$arr[10,0] = "Carla Betancourt"
$arr[10,1] = 0
$arr[10,2] = 0
$arr[10,3] = 0
$arr[10,4] = 3
$arr[10,5] = 2
$arr[10,6] = 0
$arr[11,0] = "Leah Baetcke"
$arr[11,1] = 2
$arr[11,2] = 1
$arr[11,3] = 0
$arr[11,4] = 2
$arr[11,5] = 5
$arr[11,6] = 2
$arr[12,0] = "Cason Duszak"
$arr[12,1] = 4
$arr[12,2] = 3
$arr[12,3] = 0
$arr[12,4] = 2
$arr[12,5] = 2
$arr[12,6] = 1
$arr[13,0] = "Nathan Snow"
$arr[13,1] = 4
$arr[13,2] = 0
$arr[13,3] = 0
$arr[13,4] = 2
$arr[13,5] = 3
$arr[13,6] = 1
$arr[14,0] = "Rohan Chowla"
$arr[14,1] = 17
$arr[14,2] = 6
$arr[14,3] = 1
$arr[14,4] = 2
$arr[14,5] = 0
$arr[14,6] = 0
$arr[15,0] = "Kevin Cooper"
$arr[15,1] = 8
$arr[15,2] = 4
$arr[15,3] = 0
$arr[15,4] = 2
$arr[15,5] = 1
$arr[15,6] = 0
$arr[16,0] = "Reagan Fryatt"
$arr[16,1] = 0
$arr[16,2] = 0
$arr[16,3] = 0
$arr[16,4] = 2
$arr[16,5] = 2
$arr[16,6] = 0
$arr[17,0] = "Luci Nguyen"
$arr[17,1] = 1
$arr[17,2] = 0
$arr[17,3] = 0
$arr[17,4] = 1
$arr[17,5] = 1
$arr[17,6] = 1
$arr[18,0] = "Matthew Rusten"
$arr[18,1] = 2
$arr[18,2] = 1
$arr[18,3] = 0
$arr[18,4] = 1
$arr[18,5] = 2
$arr[18,6] = 1
$arr[19,0] = "Anna Brown"
$arr[19,1] = 0
$arr[19,2] = 0
$arr[19,3] = 0
$arr[19,4] = 1
$arr[19,5] = 2
$arr[19,6] = 0
$arr[20,0] = "Piper Parker"
$arr[20,1] = 0
$arr[20,2] = 1
$arr[20,3] = 0
$arr[20,4] = 1
$arr[20,5] = 3
$arr[20,6] = 0
$arr[21,0] = "Helen Dunn"
$arr[21,1] = 1
$arr[21,2] = 0
$arr[21,3] = 0
$arr[21,4] = 0
$arr[21,5] = 1
$arr[21,6] = 1
$arr[22,0] = "Eric LastName"
$arr[22,1] = 3
$arr[22,2] = 0
$arr[22,3] = 0
$arr[22,4] = 0
$arr[22,5] = 1
$arr[22,6] = 1
$arr[23,0] = "Noah Dale"
$arr[23,1] = 1
$arr[23,2] = 0
$arr[23,3] = 0
$arr[23,4] = 0
$arr[23,5] = 3
$arr[23,6] = 1
$arr[24,0] = "Kristian Banlaoi"
$arr[24,1] = 1
$arr[24,2] = 2
$arr[24,3] = 0
$arr[24,4] = 0
$arr[24,5] = 0
$arr[24,6] = 0
$arr[25,0] = "Rose Roché"
$arr[25,1] = 3
$arr[25,2] = 6
$arr[25,3] = 0
$arr[25,4] = 0
$arr[25,5] = 0
$arr[25,6] = 0
$arr[26,0] = "Sam Tellis"
$arr[26,1] = 0
$arr[26,2] = 1
$arr[26,3] = 0
$arr[26,4] = 0
$arr[26,5] = 1
$arr[26,6] = 0
$arr[27,0] = "Cassie Deering"
$arr[27,1] = 0
$arr[27,2] = 1
$arr[27,3] = 0
$arr[27,4] = 0
$arr[27,5] = 1
$arr[27,6] = 0
$arr[28,0] = "Alex LastName"
$arr[28,1] = 1
$arr[28,2] = 1
$arr[28,3] = 0
$arr[28,4] = 0
$arr[28,5] = 1
$arr[28,6] = 0
$arr[29,0] = "Brian Tafazoli"
$arr[29,1] = 0
$arr[29,2] = 0
$arr[29,3] = 0
$arr[29,4] = 0
$arr[29,5] = 2
$arr[29,6] = 0
$arr[30,0] = "Yafu LastName"
$arr[30,1] = 0
$arr[30,2] = 0
$arr[30,3] = 0
$arr[30,4] = 0
$arr[30,5] = 2
$arr[30,6] = 0
$arr[31,0] = "Kim LastName"
$arr[31,1] = 0
$arr[31,2] = 0
$arr[31,3] = 0
$arr[31,4] = 0
$arr[31,5] = 2
$arr[31,6] = 0
$arr[32,0] = "Julie Jackson"
$arr[32,1] = 1
$arr[32,2] = 0
$arr[32,3] = 0
$arr[32,4] = 0
$arr[32,5] = 2
$arr[32,6] = 0
$arr[33,0] = "Carolyn LastName"
$arr[33,1] = 1
$arr[33,2] = 0
$arr[33,3] = 0
$arr[33,4] = 0
$arr[33,5] = 2
$arr[33,6] = 0
$arr[34,0] = "Evan Sooklal"
$arr[34,1] = 0
$arr[34,2] = 1
$arr[34,3] = 0
$arr[34,4] = 0
$arr[34,5] = 3
$arr[34,6] = 0
$arr[35,0] = "Paul Bartenfeld"
$arr[35,1] = 0
$arr[35,2] = 2
$arr[35,3] = 0
$arr[35,4] = 0
$arr[35,5] = 6
$arr[35,6] = 0

$c1 = $ws.Cells.Item(2,1)
$c2 = $ws.Cells.Item(2+$rows-1, $cols)
$rng = $ws.Range($c1, $c2)
$rng.Value = $arr

# The table grew from 30 to 37 data rows (29 to 36 players); the seven brand new
# rows (31-37) need the same formatting (bold, centered, bordered) as the rest of
# column A, which previously only existed for rows up to 30.
$fmtSrc = $ws.Range("A2")
$fmtDst = $ws.Range("A31:A37")
$fmtSrc.Copy()
$fmtDst.PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Host "Done. Range:" $rng.Address() "UsedRange:" $ws.UsedRange.Address()